$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "FirstDate" column (F) entirely.
$ws.Columns("F").Delete()

# Insert a new row for the CNY onshore commercial forward points ticker,
# directly after the existing "USDCNY 12 Month NDF Points" row (row 9).
$ws.Rows("10").Insert()
$ws.Range("A10").Value = "CNY"
$ws.Range("B10").Value = "1y"
$ws.Range("C10").Value = "fwd"
$ws.Range("D10").Value = "CCO12M BGN Curncy"
$ws.Range("E10").Value = "CNY Onsh Comm Fwd Pt 12M"

# Update the (hidden) filter-database defined name to match the new extent
# of the data (A1:E47, since a row was added and a column removed).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$E`$47"
    }
}

# Reflect the cursor position saved in the workbook at time of edit.
$ws.Range("D10").Select()
